# Update test file with empty cells
#
# - mainTimeline: clear the "-" placeholder text out of E2, E3, E4, E5,
#   E10, F10, G10 and C11 (cell formatting/style stays as-is, only the
#   stored value is removed).
# - Move the active sheet/tab selection from studyDesign to mainTimeline,
#   with the active cell on mainTimeline set to E5.

$wb = $excel.ActiveWorkbook

$timeline = $wb.Worksheets.Item("mainTimeline")

# Clear the placeholder "-" values, leaving the cell style untouched.
$timeline.Range("E2").Value = $null
$timeline.Range("E3").Value = $null
$timeline.Range("E4").Value = $null
$timeline.Range("E5").Value = $null
$timeline.Range("E10").Value = $null
$timeline.Range("F10").Value = $null
$timeline.Range("G10").Value = $null
$timeline.Range("C11").Value = $null

# Make mainTimeline the active/selected sheet and select E5 (the
# bottom-right frozen pane's active cell).
$timeline.Activate() | Out-Null
$timeline.Range("E5").Select() | Out-Null
